$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts existing rows 4-6 (and their
# values/styles) down to rows 5-7, matching the target diff.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with its data (column B stays blank).
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 3).Value = 45896.68595259259
$ws.Cells.Item(4, 4).Value = "Y2FlMmFjMDMtYmU0ZS00NzQwLTliMmMtOWM2OWUyZTA3NzAzOjU3MDE2"
